# Content_Testing_LeapFrog-games.xlsx — "Push code to review"
#
# Three content fixes on the "LeapFrog-games" sheet:
#   A5  : "(French) Football : Championnat des maths" -> "(French) Football : Khampionnat des maths"
#   B11 : "7-12 years" -> "7-15 years"
#   C22 : "$10.00" -> "$11.00" (entered with a leading apostrophe so it stays
#          text/quote-prefixed, same as the source cell's sibling price cells)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LeapFrog-games")

$ws.Range("A5").Value = "(French) Football : Khampionnat des maths"
$ws.Range("B11").Value = "7-15 years"
$ws.Range("C22").Value = "'$11.00"

# Leave the selection on the last cell touched, like a human editor would.
[void]$ws.Range("C22").Select()
